$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 102; $r -le 140; $r++) {
    $ws.Cells.Item($r, 1).Value = 0.0025
    $ws.Cells.Item($r, 2).Value = -0.25
}
